$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3858.125
$ws.Range("I40").Value = 3923.8
$ws.Range("K40").Value = 3923.8
$ws.Range("M40").Value = -3748.8
$ws.Range("H86").Value = 52657876
$ws.Range("I86").Value = 2569.3
$ws.Range("K86").Value = 2569.3
$ws.Range("M86").Value = -1446.3
$ws.Range("H88").Value = 399.5
$ws.Range("I88").Value = 399.5
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 399.5
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 6.5
$ws.Range("N88").Value = 0
$ws.Range("H89").Value = 52657876
$ws.Range("I89").Value = 2569.3
$ws.Range("K89").Value = 12846.5
$ws.Range("M89").Value = -7230.5
$ws.Range("H91").Value = 399.5
$ws.Range("I91").Value = 399.5
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 399.5
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 1004.5
$ws.Range("N91").Value = 0
$ws.Range("H99").Value = 1369299.4
$ws.Range("I99").Value = 1823257.9
$ws.Range("K99").Value = 5469773.699999999
$ws.Range("M99").Value = -5468275.699999999
$ws.Range("H100").Value = 93805.625
$ws.Range("I100").Value = 100891.4
$ws.Range("J100").Value = 81996
$ws.Range("K100").Value = 100891.4
$ws.Range("L100").Value = 81996
$ws.Range("M100").Value = -100350.4
$ws.Range("N100").Value = -83078
$ws.Range("H112").Value = 2318.0952
$ws.Range("I112").Value = 1344.3334
$ws.Range("J112").Value = 2480.389
$ws.Range("K112").Value = 4033.0002
$ws.Range("L112").Value = 7441.167
$ws.Range("M112").Value = -2925.0002
$ws.Range("N112").Value = -9657.167000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3869.3062
$ws.Range("I32").Value = 3741.5833
$ws.Range("K32").Value = 3741.5833
$ws.Range("M32").Value = -3454.5833
$ws.Range("H34").Value = 50000
$ws.Range("I34").Value = 50000
$ws.Range("K34").Value = 50000
$ws.Range("M34").Value = -49729
$ws.Range("H110").Value = 1372.1
$ws.Range("I110").Value = 810.0417
$ws.Range("J110").Value = 3620.3333
$ws.Range("K110").Value = 810.0417
$ws.Range("L110").Value = 3620.3333
$ws.Range("M110").Value = 1234.9583
$ws.Range("N110").Value = -7710.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4657.857
$ws.Range("I86").Value = 5996.9443
$ws.Range("J86").Value = 2247.5
$ws.Range("K86").Value = 5996.9443
$ws.Range("L86").Value = 2247.5
$ws.Range("M86").Value = -4873.9443
$ws.Range("N86").Value = -4493.5
$ws.Range("H89").Value = 4657.857
$ws.Range("I89").Value = 5996.9443
$ws.Range("J89").Value = 2247.5
$ws.Range("K89").Value = 29984.7215
$ws.Range("L89").Value = 11237.5
$ws.Range("M89").Value = -24368.7215
$ws.Range("N89").Value = -22469.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7827.7144
$ws.Range("I31").Value = 8659.1
$ws.Range("J31").Value = 5749.25
$ws.Range("K31").Value = 8659.1
$ws.Range("L31").Value = 5749.25
$ws.Range("M31").Value = -8364.1
$ws.Range("N31").Value = -6339.25
$ws.Range("H34").Value = 7827.7144
$ws.Range("I34").Value = 8659.1
$ws.Range("J34").Value = 5749.25
$ws.Range("K34").Value = 8659.1
$ws.Range("L34").Value = 5749.25
$ws.Range("M34").Value = -8457.1
$ws.Range("N34").Value = -6153.25
$ws.Range("H107").Value = 2930566.8
$ws.Range("I107").Value = 4454010
$ws.Range("K107").Value = 4454010
$ws.Range("M107").Value = -4452090
$ws.Range("H134").Value = 2700.7222
$ws.Range("I134").Value = 3138.75
$ws.Range("K134").Value = 9416.25
$ws.Range("M134").Value = -6881.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3833
$ws.Range("I54").Value = 750
$ws.Range("K54").Value = 2250
$ws.Range("M54").Value = -1691
$ws.Range("H103").Value = 5913.0527
$ws.Range("J103").Value = 5263.636
$ws.Range("L103").Value = 15790.908
$ws.Range("N103").Value = -17548.908
$ws.Range("H113").Value = 27622.25
$ws.Range("J113").Value = 36666.332
$ws.Range("L113").Value = 109998.996
$ws.Range("N113").Value = -114338.996
$ws.Range("H131").Value = 1629.0682
$ws.Range("J131").Value = 1657.5309
$ws.Range("L131").Value = 4972.5927
$ws.Range("N131").Value = -15052.5927

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 8000
$ws.Range("J38").Value = 8000
$ws.Range("L38").Value = 8000
$ws.Range("N38").Value = -8926
$ws.Range("H39").Value = 27220
$ws.Range("J39").Value = 26525
$ws.Range("L39").Value = 26525
$ws.Range("N39").Value = -27589
$ws.Range("H107").Value = 432.04544
$ws.Range("I107").Value = 480.27777
$ws.Range("K107").Value = 480.27777
$ws.Range("M107").Value = 1439.72223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 87900
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 87900
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H61").Value = 4555.9375
$ws.Range("I61").Value = 3244.12
$ws.Range("J61").Value = 9241
$ws.Range("K61").Value = 3244.12
$ws.Range("L61").Value = 9241
$ws.Range("M61").Value = -3042.12
$ws.Range("N61").Value = -9645
$ws.Range("H113").Value = 4555.9375
$ws.Range("I113").Value = 3244.12
$ws.Range("J113").Value = 9241
$ws.Range("K113").Value = 3244.12
$ws.Range("L113").Value = 9241
$ws.Range("M113").Value = -1074.12
$ws.Range("N113").Value = -13581
$ws.Range("H136").Value = 4918.6665
$ws.Range("I136").Value = 2891.1667
$ws.Range("K136").Value = 8673.500100000001
$ws.Range("M136").Value = -6123.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 48536.5
$ws.Range("I100").Value = 29104.777
$ws.Range("K100").Value = 58209.554
$ws.Range("M100").Value = -57668.554
$ws.Range("H122").Value = 4085.0667
$ws.Range("I122").Value = 1765.9474
$ws.Range("K122").Value = 5297.8422
$ws.Range("M122").Value = -2847.8422
$ws.Range("H123").Value = 65000
$ws.Range("J123").Value = 65000
$ws.Range("L123").Value = 65000
$ws.Range("N123").Value = -74800

Write-Host "Applied all changes"